# feat: add 2022-Q1 data
#
# Before:  Sheets = [ "2021-Q4", "总计" ]
# After:   Sheets = [ "2021-Q4", "2022-Q1", "总计" ]
#   - "2022-Q1" is a brand-new fund-holdings table (7 funds).
#   - "总计" (summary) gets a new first data row for "2022-Q1"
#     (count=7, value=4.8) while its previous "2021-Q4" row is kept,
#     now one row further down.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # "2021-Q4"   -> unchanged
$wsTotalOld = $wb.Worksheets.Item(2)   # currently "总计", physically reused as "2022-Q1"

# ---------------------------------------------------------------------
# 1) Repurpose the existing second sheet as the new "2022-Q1" detail
#    sheet and rename the (still-to-be-created) summary sheet "总计".
# ---------------------------------------------------------------------
$wsTotalOld.Name = "2022-Q1"
$ws2 = $wsTotalOld

# Clear everything so leftover "总计" cells don't linger outside the
# new table's range.
$ws2.Cells.Clear() | Out-Null

# Insert the brand new summary sheet right after "2022-Q1" (i.e. before
# where "总计" used to sit) and give it the "总计" name.
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "总计"

# ---------------------------------------------------------------------
# 2) Populate "2022-Q1" with the fund holdings table.
# ---------------------------------------------------------------------
$ws2.Range("B1").Value = "基金代码"
$ws2.Range("C1").Value = "基金名称"
$ws2.Range("D1").Value = "基金规模"
$ws2.Range("E1").Value = "股票总仓位"
$ws2.Range("F1").Value = "仓位占比"
$ws2.Range("G1").Value = "持有市值(亿元)"
$ws2.Range("H1").Value = "仓位排名"
$ws2.Range("B1:H1").Font.Bold = $true
$ws2.Range("B1:H1").HorizontalAlignment = -4108

# Numeric-looking data (fund code / scale / position / value) must stay
# TEXT, same as the "2021-Q4" sheet already does - force text format
# before assigning so Excel doesn't silently coerce them to numbers
# (which would also drop leading zeros in fund codes).
$dataCols = @("B", "C", "D", "E", "F", "G")
foreach ($col in $dataCols) {
    $ws2.Range("$($col)2:$($col)8").NumberFormat = "@"
}

$rows = @(
    @(0, "900011", "中信证券红利价值一年持有混合A", "88.99", "89.05", "2.55", "2.2692", 9),
    @(1, "900099", "中信证券红利价值一年持有混合B", "63.37", "89.05", "2.55", "1.6159", 9),
    @(2, "900089", "中信证券红利价值一年持有混合C", "22.08", "89.05", "2.55", "0.5630", 9),
    @(3, "006348", "银华盛利混合",                 "8.41",  "88.43", "2.89", "0.2430", 8),
    @(4, "013899", "上投摩根全景优势股票A",         "3.32",  "46.28", "1.96", "0.0651", 6),
    @(5, "004205", "东方支柱产业灵活配置混合",       "0.96",  "84.14", "4.28", "0.0411", 8),
    @(6, "013900", "上投摩根全景优势股票C",         "0.20",  "46.28", "1.96", "0.0039", 6)
)

$r = 2
foreach ($row in $rows) {
    $ws2.Range("A$r").Value = $row[0]
    $ws2.Range("B$r").Value = $row[1]
    $ws2.Range("C$r").Value = $row[2]
    $ws2.Range("D$r").Value = $row[3]
    $ws2.Range("E$r").Value = $row[4]
    $ws2.Range("F$r").Value = $row[5]
    $ws2.Range("G$r").Value = $row[6]
    $ws2.Range("H$r").Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 3) Populate "总计" - new 2022-Q1 row on top, the prior 2021-Q4 row
#    (previously row 2) now sits at row 3.
# ---------------------------------------------------------------------
$ws3.Range("B1").Value = "日期"
$ws3.Range("C1").Value = "持有数量(只)"
$ws3.Range("D1").Value = "持有市值(亿元)"
$ws3.Range("B1:D1").Font.Bold = $true
$ws3.Range("B1:D1").HorizontalAlignment = -4108

$ws3.Range("A2").Value = 0
$ws3.Range("B2").Value = "2022-Q1"
$ws3.Range("C2").Value = 7
$ws3.Range("D2").Value = 4.8

$ws3.Range("A3").Value = 1
$ws3.Range("B3").Value = "2021-Q4"
$ws3.Range("C3").Value = 2
$ws3.Range("D3").Value = 0

$ws1.Range("A1").Select() | Out-Null
